# Applies the "Entrega 4" justification rows + related formatting tweaks
# described by the commit "cambio justificaciones y diagrama".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: new section header "Entrega 4" (reuses the existing style of
#     the surrounding blank rows, s="5") ---------------------------------
$ws.Range("A21").Value = "Entrega 4"

# --- Row 22: Vinculador decision -----------------------------------------
$ws.Range("A22").Value = "Objetos"
$ws.Range("B22").Value = "Vinculador"
$ws.Range("C22").Value = "Encapsulamos el proceso del validador en una clase que se encarga de recibir la entidad y los criterios con los que quiere vincular"
$ws.Range("D22").Value = "Hacer un Strategy y una clase para cada metodo de validacion"

# --- Row 23: Item / ItemEgreso / ItemPresupuesto decision -----------------
$ws.Range("A23").Value = "Objetos"
$ws.Range("B23").Value = "Item, ItemEgreso, ItemPresupuesto"
$ws.Range("C23").Value = "Sacamos la clase Item y Separamos al item del egreso del item del ingreso para tratarlos como objetos diferentes"
$ws.Range("D23").Value = "Hacer que ItemEgreso e ItemPresupuesto hereden de la misma clase item, no pudimos hacerlo asi por limitaciones del ORM"

# E23 gets a new wrap-text + underlined-font style (matches the new
# cellXfs/fonts entries added to styles.xml).
$ws.Range("E23").Font.Underline = $true

# --- Reflect the author's final cursor position ---------------------------
$ws.Range("E23").Select()
